# Adds more tests for database
# - Update row 2: lower the average mark and correct the teacher's patronymic
# - Remove the old row 3 (Winter/2019, Shibeko Viktor Nikolaevich) and shift
#   the former row 4 (Winter/2020, Chaikovski Petr Ilich) up into row 3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.5
$ws.Range("C2").Value = "Mocart Amadey Batkovich"

$ws.Rows(3).Delete()
